$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{ A = "HarveyCanalNorth"; B = "Harvey Canal Sector Gates North / Prot Side nr Lapalco Blvd"; C = 2 },
    @{ A = "HarveyCanalBoom";  B = "Harvey Canal at Boomtown Casion"; C = 0 },
    @{ A = "BayouBienv";       B = "Bayou Bienvenue Floodgate"; C = 3 },
    @{ A = "BaraPass";         B = "Barataria Pass at Grand Isle"; C = 0 },
    @{ A = "FreshCanal";       B = "Freshwater Canal at Freshwater Bayou Lock South"; C = 0 },
    @{ A = "CalcRiv";          B = "Calcasieu River at Cameron"; C = 0 }
)

$startRow = 17

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
}

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $data.B
}

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 3).NumberFormat = "0.0"
}

# Columns A and B now hold longer text (IDs + new long station names) -
# widen them to fit the new content, like a user re-running "AutoFit Column Width".
$ws.Columns.Item(1).ColumnWidth = 15.41796875
$ws.Columns.Item(2).ColumnWidth = 47.9453125

# Move the active selection to the first empty row below the new data,
# matching where the user's cursor ended up after entering the rows.
$ws.Range("C23").Select() | Out-Null

